$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at row 34 (pushes existing rows 34..60 down to 35..61) ---
$ws.Rows("34").Insert()

# --- 2. Copy the formatting of the existing DNP template row (R76 / "RESISTOR" DNP row,
#        which used to be row 44 and is now row 45 after the insert above) onto the new
#        row 34, restricted to columns A:L so we don't blow out the whole row. ---
$ws.Range("A45:L45").Copy()
$ws.Range("A34:L34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Populate the new row 34 with the R14 / DNP details. "R14" must be written
#        (added as a new shared string) before the trimmed reference text below so
#        the shared-string table ends up ordered the same way Excel produced it. ---
$ws.Range("A34").Value = 30
$ws.Range("B34").Value = "R14"
$ws.Range("C34").Value = "DNP"
$ws.Range("D34").Value = "R0402"
$ws.Range("E34").Value = "Rohm"
$ws.Range("F34").Value = "TRR01MZPJ000"
$ws.Range("G34").Value = "Digi-Key"
$ws.Range("H34").Value = "RHM0.0BICT-ND"
$ws.Range("I34").Value = "RES 0.0 OHM 1/16W 0402 SMD"
$ws.Range("J34").Value = 1
$ws.Range("K34").Value = 0.181
$ws.Range("L34").Formula = "=J34*K34"

# --- 4. Remove "R14" from the reference list in row 33 (it now gets its own DNP line). ---
$ws.Range("B33").Value = "R4, R11, R12, R13, R15, R16, R17, R18, R49, R74, R75"

# --- 5. Restore the sheet view scroll/selection state seen after the edit ---
$ws.Range("E22").Select()
$excel.ActiveWindow.ScrollRow = 7
